$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value2 = $value
}

function Clear-Cell($ws, $row, $col) {
    $ws.Cells.Item($row, $col).Value2 = $null
}

# ---------------- ALC (sheet1) ----------------
$ws = $wb.Worksheets.Item("ALC")

# Row 96
Set-Cell $ws 96 8 2763.6538
Set-Cell $ws 96 9 3372.111
Set-Cell $ws 96 11 10116.333
Set-Cell $ws 96 13 -8743.332999999999

# Row 99
Set-Cell $ws 99 8 1158
Set-Cell $ws 99 9 737.25
Set-Cell $ws 99 10 1999.5
Set-Cell $ws 99 11 2211.75
Set-Cell $ws 99 12 5998.5
Set-Cell $ws 99 13 -713.75
Set-Cell $ws 99 14 -8994.5

# Row 111
Set-Cell $ws 111 8 9994.143
Set-Cell $ws 111 10 4327.5
Set-Cell $ws 111 12 12982.5
Set-Cell $ws 111 14 -19116.5

# Row 137
Set-Cell $ws 137 8 1000
Set-Cell $ws 137 9 1000
Set-Cell $ws 137 10 0
Set-Cell $ws 137 11 3000
Set-Cell $ws 137 12 0
Set-Cell $ws 137 13 -450
Clear-Cell $ws 137 14

# ---------------- ARM (sheet2) ----------------
$ws = $wb.Worksheets.Item("ARM")

# Row 97
Set-Cell $ws 97 8 857.6667
Set-Cell $ws 97 9 857.3333
Set-Cell $ws 97 10 859.3333
Set-Cell $ws 97 11 857.3333
Set-Cell $ws 97 12 859.3333
Set-Cell $ws 97 13 -361.3333
Set-Cell $ws 97 14 -1851.3333

# Row 113
Set-Cell $ws 113 8 149997
Set-Cell $ws 113 10 149997
Set-Cell $ws 113 12 149997
Set-Cell $ws 113 14 -158675

# ---------------- BSM (sheet3) ----------------
$ws = $wb.Worksheets.Item("BSM")

# Row 5
Set-Cell $ws 5 8 505.66666
Set-Cell $ws 5 9 206.8
Set-Cell $ws 5 10 2000
Set-Cell $ws 5 11 206.8
Set-Cell $ws 5 12 2000
Set-Cell $ws 5 13 -93.80000000000001
Set-Cell $ws 5 14 -2226

# ---------------- CRP (sheet4) ----------------
$ws = $wb.Worksheets.Item("CRP")

# Row 16
Set-Cell $ws 16 8 5253.6
Set-Cell $ws 16 10 6299.6
Set-Cell $ws 16 12 6299.6
Set-Cell $ws 16 14 -6873.6

# Row 31
Set-Cell $ws 31 8 1303
Set-Cell $ws 31 9 1303
Set-Cell $ws 31 10 0
Set-Cell $ws 31 11 1303
Set-Cell $ws 31 12 0
Set-Cell $ws 31 13 -1008
Clear-Cell $ws 31 14

# Row 34
Set-Cell $ws 34 8 1303
Set-Cell $ws 34 9 1303
Set-Cell $ws 34 10 0
Set-Cell $ws 34 11 1303
Set-Cell $ws 34 12 0
Set-Cell $ws 34 13 -1101
Clear-Cell $ws 34 14

# Row 113
Set-Cell $ws 113 8 5253.6
Set-Cell $ws 113 10 6299.6
Set-Cell $ws 113 12 6299.6
Set-Cell $ws 113 14 -10639.6

# ---------------- CUL (sheet5) ----------------
$ws = $wb.Worksheets.Item("CUL")

# Row 4
Set-Cell $ws 4 8 9253202
Set-Cell $ws 4 9 10910820
Set-Cell $ws 4 10 3175270
Set-Cell $ws 4 11 32732460
Set-Cell $ws 4 12 9525810
Set-Cell $ws 4 13 -32732348
Set-Cell $ws 4 14 -9526034

# Row 80
Set-Cell $ws 80 8 2049.4
Set-Cell $ws 80 10 2081.6667
Set-Cell $ws 80 12 6245.000100000001
Set-Cell $ws 80 14 -8117.000100000001

# Row 83
Set-Cell $ws 83 8 2049.4
Set-Cell $ws 83 10 2081.6667
Set-Cell $ws 83 12 18735.0003
Set-Cell $ws 83 14 -28095.0003

# Row 86
Set-Cell $ws 86 8 0
Set-Cell $ws 86 9 0
Set-Cell $ws 86 10 0
Set-Cell $ws 86 11 0
Set-Cell $ws 86 12 0
Clear-Cell $ws 86 13
Clear-Cell $ws 86 14

# Row 88
Set-Cell $ws 88 8 17999.4
Set-Cell $ws 88 10 17999.4
Set-Cell $ws 88 12 53998.2
Set-Cell $ws 88 14 -54854.2

# Row 89
Set-Cell $ws 89 8 0
Set-Cell $ws 89 9 0
Set-Cell $ws 89 10 0
Set-Cell $ws 89 11 0
Set-Cell $ws 89 12 0
Clear-Cell $ws 89 13
Clear-Cell $ws 89 14

# Row 91
Set-Cell $ws 91 8 17999.4
Set-Cell $ws 91 10 17999.4
Set-Cell $ws 91 12 53998.2
Set-Cell $ws 91 14 -56962.2

# Row 95
Set-Cell $ws 95 8 8875.5
Set-Cell $ws 95 10 8875.5
Set-Cell $ws 95 12 26626.5
Set-Cell $ws 95 14 -30744.5

# Row 119
Set-Cell $ws 119 8 0
Set-Cell $ws 119 9 0
Set-Cell $ws 119 11 0
Clear-Cell $ws 119 13

# Row 120
Set-Cell $ws 120 8 0
Set-Cell $ws 120 9 0
Set-Cell $ws 120 11 0
Clear-Cell $ws 120 13

# Row 123
Set-Cell $ws 123 8 2733
Set-Cell $ws 123 9 2733
Set-Cell $ws 123 11 8199
Set-Cell $ws 123 13 -5749

# ---------------- GSM (sheet6) ----------------
$ws = $wb.Worksheets.Item("GSM")

# Row 18
Set-Cell $ws 18 8 55000
Set-Cell $ws 18 9 55000
Set-Cell $ws 18 11 55000
Set-Cell $ws 18 13 -54707

# Row 33
Set-Cell $ws 33 8 30025000
Set-Cell $ws 33 10 30025000
Set-Cell $ws 33 12 30025000
Set-Cell $ws 33 14 -30025504

# Row 80
Set-Cell $ws 80 8 3153.3333
Set-Cell $ws 80 9 2935
Set-Cell $ws 80 10 4900
Set-Cell $ws 80 11 2935
Set-Cell $ws 80 12 4900
Set-Cell $ws 80 13 -1937
Set-Cell $ws 80 14 -6896

# Row 83
Set-Cell $ws 83 8 3153.3333
Set-Cell $ws 83 9 2935
Set-Cell $ws 83 10 4900
Set-Cell $ws 83 11 14675
Set-Cell $ws 83 12 24500
Set-Cell $ws 83 13 -9683
Set-Cell $ws 83 14 -34484

# Row 136
Set-Cell $ws 136 8 26500
Set-Cell $ws 136 10 26500
Set-Cell $ws 136 12 79500
Set-Cell $ws 136 14 -84600

# ---------------- LTW (sheet7) ----------------
$ws = $wb.Worksheets.Item("LTW")

# Row 10
Set-Cell $ws 10 8 752499.75
Set-Cell $ws 10 9 1002999.7
Set-Cell $ws 10 10 1000
Set-Cell $ws 10 11 1002999.7
Set-Cell $ws 10 12 1000
Set-Cell $ws 10 13 -1002859.7
Set-Cell $ws 10 14 -1280

# Row 136
Set-Cell $ws 136 8 2066.818
Set-Cell $ws 136 10 2000
Set-Cell $ws 136 12 6000
Set-Cell $ws 136 14 -11100

# ---------------- WVR (sheet8) ----------------
$ws = $wb.Worksheets.Item("WVR")

# Row 7
Set-Cell $ws 7 8 2500
Set-Cell $ws 7 9 500
Set-Cell $ws 7 10 4500
Set-Cell $ws 7 11 500
Set-Cell $ws 7 12 4500
Set-Cell $ws 7 13 -387
Set-Cell $ws 7 14 -4726

# Row 132
Set-Cell $ws 132 8 846.5833
Set-Cell $ws 132 9 941.1111
Set-Cell $ws 132 11 2823.3333
Set-Cell $ws 132 13 -293.3332999999998
